$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.136730802769335
$ws.Range("C2").Value = 0.1710915350730886
$ws.Range("D2").Value = 0.05164015962813551
$ws.Range("E2").Value = 0.1237518460948479
$ws.Range("F2").Value = 3.297848894863449
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.226129125654424
$ws.Range("K2").Value = 1.229727015956598
$ws.Range("N2").Value = 3.111090195120568

$ws.Range("B3").Value = 1.090893794756511
$ws.Range("C3").Value = 0.1629003829942803
$ws.Range("D3").Value = 0.05030122963908212
$ws.Range("E3").Value = 0.1211336060323056
$ws.Range("F3").Value = 3.27336112416242
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.2216855784172935
$ws.Range("K3").Value = 1.178190574632595
$ws.Range("N3").Value = 3.116910752226474

$ws.Range("B4").Value = 1.063452048856703
$ws.Range("C4").Value = 0.1579826333693006
$ws.Range("D4").Value = 0.04951307680556027
$ws.Range("E4").Value = 0.1195983299881931
$ws.Range("F4").Value = 3.259770387011159
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.2190838162325051
$ws.Range("K4").Value = 1.147316752513689
$ws.Range("N4").Value = 3.121222168917527

$ws.Range("B5").Value = 1.052445612451834
$ws.Range("C5").Value = 0.1560065759220208
$ws.Range("D5").Value = 0.04920046178716575
$ws.Range("E5").Value = 0.1189908644581266
$ws.Range("F5").Value = 3.25459483534442
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.2180553356912327
$ws.Range("K5").Value = 1.13492855028926
$ws.Range("N5").Value = 3.123164253034929

$ws.Range("B6").Value = 1.050628641023849
$ws.Range("C6").Value = 0.1556801389637457
$ws.Range("D6").Value = 0.04914907027839632
$ws.Range("E6").Value = 0.1188910923315092
$ws.Range("F6").Value = 3.253757334788006
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.2178864732843309
$ws.Range("K6").Value = 1.132883151052397
$ws.Range("N6").Value = 3.123497907706479

$ws.Range("B7").Value = 1.063302898798753
$ws.Range("C7").Value = 0.1579558704803077
$ws.Range("D7").Value = 0.04950882606201645
$ws.Range("E7").Value = 0.1195900639395902
$ws.Range("F7").Value = 3.259699119536137
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.2190698172906878
$ws.Range("K7").Value = 1.147148899472512
$ws.Range("N7").Value = 3.121247611338802

$ws.Range("B8").Value = 1.120780216510354
$ws.Range("C8").Value = 0.1682439829780265
$ws.Range("D8").Value = 0.05117146351055624
$ws.Range("E8").Value = 0.1228340523967653
$ws.Range("F8").Value = 3.289105292169026
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.2245706690114986
$ws.Range("K8").Value = 1.211797204968917
$ws.Range("N8").Value = 3.112943843857465

$ws.Range("B9").Value = 1.239088032940288
$ws.Range("C9").Value = 0.1893111294020855
$ws.Range("D9").Value = 0.0547004706563925
$ws.Range("E9").Value = 0.1297707430805914
$ws.Range("F9").Value = 3.358265995471655
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2363665225199583
$ws.Range("K9").Value = 1.344708142913817
$ws.Range("N9").Value = 3.10252808668659

$ws.Range("B10").Value = 1.329461333774304
$ws.Range("C10").Value = 0.205343685082056
$ws.Range("D10").Value = 0.05745622924138871
$ws.Range("E10").Value = 0.1352203846030164
$ws.Range("F10").Value = 3.41613906527553
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2456550031429572
$ws.Range("K10").Value = 1.446149712523209
$ws.Range("N10").Value = 3.098476087102441

$ws.Range("B11").Value = 1.371333998530247
$ws.Range("C11").Value = 0.2127601558477181
$ws.Range("D11").Value = 0.0587451651277533
$ws.Range("E11").Value = 0.1377768884715493
$ws.Range("F11").Value = 3.444012401755089
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2500172867450203
$ws.Range("K11").Value = 1.493133486145297
$ws.Range("N11").Value = 3.097419738053659

$ws.Range("B12").Value = 1.387300156466438
$ws.Range("C12").Value = 0.215586455195421
$ws.Range("D12").Value = 0.05923831415212533
$ws.Range("E12").Value = 0.1387561407712994
$ws.Range("F12").Value = 3.454790530479585
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2516889667091675
$ws.Range("K12").Value = 1.511046184117447
$ws.Range("N12").Value = 3.097133297988236

$ws.Range("B13").Value = 1.38385666485965
$ws.Range("C13").Value = 0.2149769654069189
$ws.Range("D13").Value = 0.05913188113790113
$ws.Range("E13").Value = 0.1385447445302077
$ws.Range("F13").Value = 3.452459333764835
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2513280594320548
$ws.Range("K13").Value = 1.507182975654189
$ws.Range("N13").Value = 3.097189929922308

$ws.Range("B14").Value = 1.37264533928726
$ws.Range("C14").Value = 0.2129923188485918
$ws.Range("D14").Value = 0.05878563557532601
$ws.Range("E14").Value = 0.1378572283441031
$ws.Range("F14").Value = 3.444894649531932
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2501544199587187
$ws.Range("K14").Value = 1.494604747223576
$ws.Range("N14").Value = 3.097393893805901

$ws.Range("B15").Value = 1.36579240147114
$ws.Range("C15").Value = 0.2117789930776439
$ws.Range("D15").Value = 0.05857420806309932
$ws.Range("E15").Value = 0.1374375589100651
$ws.Range("F15").Value = 3.44029013749585
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2494381102727914
$ws.Range("K15").Value = 1.486915988033218
$ws.Range("N15").Value = 3.097533631061381

$ws.Range("B16").Value = 1.326740211895924
$ws.Range("C16").Value = 0.2048614910086997
$ws.Range("D16").Value = 0.05737270314649834
$ws.Range("E16").Value = 0.1350548713420565
$ws.Range("F16").Value = 3.414348658304988
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2453726802944942
$ws.Range("K16").Value = 1.443096110118404
$ws.Range("N16").Value = 3.098560992498818

$ws.Range("B17").Value = 1.302978251960781
$ws.Range("C17").Value = 0.2006494685108748
$ws.Range("D17").Value = 0.05664465043857803
$ws.Range("E17").Value = 0.1336130198966217
$ws.Range("F17").Value = 3.398831015317938
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2429137999590552
$ws.Range("K17").Value = 1.41642891106946
$ws.Range("N17").Value = 3.099393091108681

$ws.Range("B18").Value = 1.289382642732278
$ws.Range("C18").Value = 0.1982384110322926
$ws.Range("D18").Value = 0.05622922031631816
$ws.Range("E18").Value = 0.1327909920583892
$ws.Range("F18").Value = 3.390051206649844
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2415123978554732
$ws.Range("K18").Value = 1.401169439868312
$ws.Range("N18").Value = 3.099945736105937

$ws.Range("B19").Value = 1.284791695548165
$ws.Range("C19").Value = 0.1974240544748511
$ws.Range("D19").Value = 0.05608913486300793
$ws.Range("E19").Value = 0.132513918362001
$ws.Range("F19").Value = 3.387103487831752
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2410401167382048
$ws.Range("K19").Value = 1.396016362474001
$ws.Range("N19").Value = 3.100145556679848

$ws.Range("B20").Value = 1.305500335665442
$ws.Range("C20").Value = 0.2010966454099332
$ws.Range("D20").Value = 0.05672180877702715
$ws.Range("E20").Value = 0.1337657530927459
$ws.Range("F20").Value = 3.400467828184901
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2431742181983765
$ws.Range("K20").Value = 1.419259521004562
$ws.Range("N20").Value = 3.099296846573353

$ws.Range("B21").Value = 1.375935391781525
$ws.Range("C21").Value = 0.2135747723572194
$ws.Range("D21").Value = 0.05888719922650409
$ws.Range("E21").Value = 0.1380588654390138
$ws.Range("F21").Value = 3.447110520079974
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2504986088253958
$ws.Range("K21").Value = 1.49829599061988
$ws.Range("N21").Value = 3.09733089888023

$ws.Range("B22").Value = 1.422609481292113
$ws.Range("C22").Value = 0.2218340132267826
$ws.Range("D22").Value = 0.0603318726437152
$ws.Range("E22").Value = 0.1409297223353079
$ws.Range("F22").Value = 3.47889495654897
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2554008451648997
$ws.Range("K22").Value = 1.550656211369471
$ws.Range("N22").Value = 3.096708246632645

$ws.Range("B23").Value = 1.397639881090811
$ws.Range("C23").Value = 0.2174163356418433
$ws.Range("D23").Value = 0.05955813464385074
$ws.Range("E23").Value = 0.1393915305878863
$ws.Range("F23").Value = 3.461811750428041
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2527738468379397
$ws.Range("K23").Value = 1.522645860179125
$ws.Range("N23").Value = 3.096979836233757

$ws.Range("B24").Value = 1.304359897760889
$ws.Range("C24").Value = 0.2008944440694052
$ws.Range("D24").Value = 0.05668691571929685
$ws.Range("E24").Value = 0.1336966808883346
$ws.Range("F24").Value = 3.399727384403633
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2430564449922628
$ws.Range("K24").Value = 1.41797957847146
$ws.Range("N24").Value = 3.099340127430338

$ws.Range("B25").Value = 1.206479383211672
$ws.Range("C25").Value = 0.1835153039756676
$ws.Range("D25").Value = 0.05371710532924112
$ws.Range("E25").Value = 0.1278323382273321
$ws.Range("F25").Value = 3.338320128606483
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.2330667112360061
$ws.Range("K25").Value = 1.308090071279821
$ws.Range("N25").Value = 3.104715246303115
